$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between row 16 and row 17
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 56940

$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 7592
